$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data / recalculated mean
$ws.Range("F2").Value = -7
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -2
$ws.Range("F7").Value = 4
$ws.Range("F10").Value = -2
$ws.Range("F14").Value = 5
$ws.Range("F16").Value = 4
$ws.Range("F17").Value = -3
$ws.Range("F18").Value = -2
